# Updated cryptos list on Fri Oct 18 19:42:11 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row.
# A leading "'" forces Excel to keep a numeric-looking price as literal
# text (matching the original inlineStr cell type) instead of coercing
# it to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.690.26"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "2.652.91"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'599.99"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").Value = "'154.97"
$ws.Range("E6").Value = "  +3.83%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.546"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "2.652.22"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("D10").Value = "'0.137"
$ws.Range("E10").Value = "  +11.84%  "
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "'5.23"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "'0.348"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "'27.94"
$ws.Range("E14").Value = "  +3.01%  "
$ws.Range("D15").Value = "'0.0000188"
$ws.Range("E15").Value = "  +6.05%  "
$ws.Range("D16").Value = "3.138.59"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "68.618.77"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("D18").Value = "2.658.44"
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("D19").Value = "'11.42"
$ws.Range("E19").Value = "  +4.06%  "
$ws.Range("D20").Value = "'366.77"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").Value = "'7.43"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'4.87"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  +4.80%  "
$ws.Range("D25").Value = "'72.65"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'9.98"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("D28").Value = "'0.0000105"
$ws.Range("E28").Value = "  +7.92%  "
$ws.Range("E29").Value = "  +2.44%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'574.11"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").Value = "'1.42"
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("D33").Value = "'7.97"
$ws.Range("E33").Value = "  +4.93%  "
$ws.Range("D34").Value = "'1.85"
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("D35").Value = "'0.130"
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "'1.54"
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("D38").Value = "'158.73"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D39").Value = "'1.93"
$ws.Range("E39").Value = "  +5.29%  "
$ws.Range("D40").Value = "'19.25"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").Value = "'5.39"
$ws.Range("E41").Value = "  +4.02%  "
$ws.Range("D42").Value = "'0.368"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("D43").Value = "'2.65"
$ws.Range("E43").Value = "  +6.93%  "
$ws.Range("D44").Value = "'17.75"
$ws.Range("E44").Value = "  +4.74%  "
$ws.Range("D45").Value = "0.0₆0320"
$ws.Range("E45").Value = "  +12.75%  "
$ws.Range("D46").Value = "'40.64"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").Value = "'156.63"
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("D49").Value = "'3.73"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("D51").Value = "'21.95"
$ws.Range("E51").Value = "  +3.33%  "
